# Apply the edits described in the commit: update the "test_suite" sheet so
# that the login test row uses "login"/"Y" and the CreateAccountTest row is
# removed entirely (the corresponding data now lives only on the
# CreateAccountTest sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_suite")

# Update row 2: logintest/N -> login/Y
$ws.Cells.Item(2, 1).Value = "login"
$ws.Cells.Item(2, 2).Value = "Y"

# Remove row 3 (CreateAccountTest / Y) entirely.
$ws.Rows.Item(3).Delete()

# Restore the originally-selected cell on this sheet.
$ws.Range("D5").Select()
